# Append simulation rows s11..s15 to the manifest worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append (subject_id, file_name, feedback_1_id, feedback_1_x, feedback_1_y,
# feedback_1_toleranceA, feedback_1_toleranceB, feedback_1_theta, minor_to_major_ratio)
$newRows = @(
    @("s11", "s11_IMG_3178.jpeg", "meltpatch", "1721", "2538", "104", "52", "51", "2"),
    @("s12", "s12_IMG_3180.jpeg", "meltpatch", "79",   "1405", "104", "52", "147", "2"),
    @("s13", "s13_IMG_3176.jpeg", "meltpatch", "1155", "2293", "104", "52", "108", "2"),
    @("s14", "s14_IMG_3179.jpeg", "meltpatch", "635",  "2056", "104", "52", "84",  "2"),
    @("s15", "s15_IMG_3174.jpeg", "meltpatch", "2777", "1232", "104", "52", "11",  "2")
)

$startRow = 12
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowNum = $startRow + $i
    $rowData = $newRows[$i]
    for ($col = 1; $col -le $rowData.Count; $col++) {
        $value = $rowData[$col - 1]
        $cell = $ws.Cells.Item($rowNum, $col)
        if ($col -ge 4) {
            # Columns D..I hold purely numeric-looking text (e.g. "1721"). Excel
            # would otherwise auto-convert these to numbers, but the sheet stores
            # every cell as text, so force the text format before assigning.
            $cell.NumberFormat = "@"
        }
        $cell.Value = $value
    }
}

$wb.Save()
